# Edit script for B0CSJYJRKD_po_data.xlsx
# - Rename "Requested quantity" header on "Weekly Quantity" sheet to "Weekly_PO_Qty"
# - Rename "Requested quantity" header on "Monthly Trend" sheet to "Monthly_PO_Qty"
# - Add a new "PO Forecast" sheet (third sheet) with ds/PO_Forecast/yhat_lower/yhat_upper columns

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename header cells
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the bold/centered/bordered header style from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the date-time number format from the Weekly Quantity sheet's date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A37").PasteSpecial(-4122)

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45312.99999999999, 109, -25.22563873156716, 243.8570249987102),
    @(45319.99999999999, 109, -30.11607340786951, 252.3513012732099),
    @(45340.99999999999, 109, -21.93844892560523, 250.3002281672516),
    @(45347.99999999999, 109, -23.70949343650914, 251.2737917969661),
    @(45410.99999999999, 109, -29.85032916059425, 250.9343561198435),
    @(45424.99999999999, 109, -28.99342454607557, 239.0859985781612),
    @(45431.99999999999, 109, -35.16149491400288, 242.2289024943926),
    @(45445.99999999999, 109, -34.57466899164118, 250.5907523090019),
    @(45452.99999999999, 109, -32.1857031538178, 245.5043278590333),
    @(45466.99999999999, 109, -34.5710882324387, 250.7507091172172),
    @(45473.99999999999, 109, -26.52098203853965, 244.9597563404804),
    @(45480.99999999999, 109, -35.15875189916346, 243.6230281445864),
    @(45487.99999999999, 109, -22.69854492279514, 247.6112238889181),
    @(45494.99999999999, 109, -35.77453030178589, 231.9659570096558),
    @(45501.99999999999, 109, -39.29431225312685, 245.494812088982),
    @(45508.99999999999, 109, -38.96202078818314, 249.2813911520444),
    @(45515.99999999999, 109, -46.60554640719557, 256.6907419040899),
    @(45522.99999999999, 109, -29.98142188986229, 251.1040686391097),
    @(45536.99999999999, 109, -29.48142784866822, 256.918555746498),
    @(45543.99999999999, 109, -32.87340733705599, 248.9689496244465),
    @(45550.99999999999, 109, -33.39193414953618, 252.3465813238842),
    @(45557.99999999999, 109, -30.94876690034762, 243.2341447741742),
    @(45564.99999999999, 109, -24.02888754333092, 258.1765743470375),
    @(45571.99999999999, 109, -23.77973878844077, 255.2814454999216),
    @(45578.99999999999, 109, -29.15431789917062, 243.1660047897002),
    @(45585.99999999999, 109, -20.9084928164606, 251.6956043051875),
    @(45592.99999999999, 109, -19.94249705830782, 244.6827829696033),
    @(45613.99999999999, 109, -38.26938141926184, 245.9123938301348),
    @(45620.99999999999, 109, -30.58434904481534, 243.9333757531647),
    @(45627.99999999999, 109, -30.07381821223786, 251.8897198301706),
    @(45634.99999999999, 109, -34.59312990036802, 237.9304646106509),
    @(45641.99999999999, 109, -30.51372791617973, 242.8537600405748),
    @(45648.99999999999, 109, -23.45772408219868, 238.0034495813442),
    @(45655.99999999999, 109, -29.13951156336996, 245.7586260358383),
    @(45662.99999999999, 109, -38.29604935274232, 248.0131661029697),
    @(45669.99999999999, 109, -31.32543924596219, 252.620959627349)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
